{"js": "// Update the date line and the two-digit \u00f7 one-digit division problems\n// to the new \"output generated at c986bee\" values.\n\nconst body = context.document.body;\n\nasync function replaceAll(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// Header date\nawait replaceAll(\"2024-10-25 Friday\", \"2024-10-26 Saturday\");\n\n// Simple 1:1 unique replacements\nawait replaceAll(\"12\u00f79=\", \"78\u00f79=\");\nawait replaceAll(\"31\u00f78=\", \"41\u00f75=\");\nawait replaceAll(\"65\u00f75=\", \"55\u00f74=\");\nawait replaceAll(\"78\u00f72=\", \"85\u00f79=\");\nawait replaceAll(\"43\u00f77=\", \"22\u00f77=\");\nawait replaceAll(\"80\u00f78=\", \"59\u00f79=\");\nawait replaceAll(\"83\u00f78=\", \"84\u00f79=\");\nawait replaceAll(\"85\u00f75=\", \"91\u00f73=\");\nawait replaceAll(\"58\u00f79=\", \"33\u00f73=\");\nawait replaceAll(\"97\u00f77=\", \"58\u00f77=\");\nawait replaceAll(\"45\u00f78=\", \"25\u00f75=\");\nawait replaceAll(\"18\u00f78=\", \"49\u00f78=\");\nawait replaceAll(\"61\u00f73=\", \"85\u00f78=\");\nawait replaceAll(\"35\u00f77=\", \"65\u00f74=\");\nawait replaceAll(\"56\u00f74=\", \"86\u00f75=\");\nawait replaceAll(\"92\u00f72=\", \"78\u00f78=\");\nawait replaceAll(\"25\u00f74=\", \"50\u00f79=\");\nawait replaceAll(\"88\u00f77=\", \"79\u00f76=\");\nawait replaceAll(\"47\u00f75=\", \"61\u00f78=\");\nawait replaceAll(\"91\u00f72=\", \"23\u00f76=\");\nawait replaceAll(\"75\u00f77=\", \"44\u00f76=\");\nawait replaceAll(\"21\u00f74=\", \"59\u00f74=\");\nawait replaceAll(\"86\u00f77=\", \"11\u00f72=\");\n\n// \"21\u00f77=\" appears twice (two different cells) and maps to two different\n// results, so handle both occurrences positionally, in document order.\n{\n  const results = body.search(\"21\u00f77=\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  const replacements = [\"95\u00f79=\", \"95\u00f76=\"];\n  for (let i = 0; i < results.items.length && i < replacements.length; i++) {\n    results.items[i].insertText(replacements[i], Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the two-digit \u00f7 one-digit division problems\n# to the new \"output generated at c986bee\" values.\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($oldText, $newText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n# Header date\nReplace-AllText \"2024-10-25 Friday\" \"2024-10-26 Saturday\"\n\n# Simple 1:1 unique replacements\nReplace-AllText \"12\u00f79=\" \"78\u00f79=\"\nReplace-AllText \"31\u00f78=\" \"41\u00f75=\"\nReplace-AllText \"65\u00f75=\" \"55\u00f74=\"\nReplace-AllText \"78\u00f72=\" \"85\u00f79=\"\nReplace-AllText \"43\u00f77=\" \"22\u00f77=\"\nReplace-AllText \"80\u00f78=\" \"59\u00f79=\"\nReplace-AllText \"83\u00f78=\" \"84\u00f79=\"\nReplace-AllText \"85\u00f75=\" \"91\u00f73=\"\nReplace-AllText \"58\u00f79=\" \"33\u00f73=\"\nReplace-AllText \"97\u00f77=\" \"58\u00f77=\"\nReplace-AllText \"45\u00f78=\" \"25\u00f75=\"\nReplace-AllText \"18\u00f78=\" \"49\u00f78=\"\nReplace-AllText \"61\u00f73=\" \"85\u00f78=\"\nReplace-AllText \"35\u00f77=\" \"65\u00f74=\"\nReplace-AllText \"56\u00f74=\" \"86\u00f75=\"\nReplace-AllText \"92\u00f72=\" \"78\u00f78=\"\nReplace-AllText \"25\u00f74=\" \"50\u00f79=\"\nReplace-AllText \"88\u00f77=\" \"79\u00f76=\"\nReplace-AllText \"47\u00f75=\" \"61\u00f78=\"\nReplace-AllText \"91\u00f72=\" \"23\u00f76=\"\nReplace-AllText \"75\u00f77=\" \"44\u00f76=\"\nReplace-AllText \"21\u00f74=\" \"59\u00f74=\"\nReplace-AllText \"86\u00f77=\" \"11\u00f72=\"\n\n# \"21\u00f77=\" appears twice (two different cells) and maps to two different\n# results, so walk the occurrences in document order and replace each in turn.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$targets = @(\"95\u00f79=\", \"95\u00f76=\")\n$i = 0\nwhile ($i -lt $targets.Length -and $rng.Find.Execute(\"21\u00f77=\")) {\n    $rng.Text = $targets[$i]\n    $rng.Collapse(0)\n    $i = $i + 1\n}\n"}
